$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPM")

# Clear the "Temps Previst" (expected time) raw input for tasks S1-S4 (rows 52-55, col J)
$ws.Range("J52:J55").ClearContents()
